$d = $word.ActiveDocument

# --- Fix the "empty note" bug: footnote 28 was an orphaned/near-empty
# note ("a") that really belongs appended to the end of footnote 27's
# text. Remove footnote 28 entirely and fold its content onto 27. ---

# Footnote 28 is the 8th footnote in the document (last one) - delete it
# outright, it contributes nothing once merged into 27.
$fn28 = $d.Footnotes.Item(8)
$fn28.Delete()

# Footnote 27 (now the 7th / last remaining footnote) needs "a" appended
# to its note text. Rebuild it in place at the same reference point so
# the run layout matches a normally-typed note (mark + space + text).
$fn27 = $d.Footnotes.Item(7)
$refRange = $fn27.Reference
$insertionPoint = $d.Range($refRange.Start, $refRange.End)
$noteText = $fn27.Range.Text
$fn27.Delete()

$d.Footnotes.Add($insertionPoint, [Type]::Missing, ($noteText + "a")) | Out-Null
